$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 304.33334
$ws.Range("I33").Value = 311.1875
$ws.Range("J33").Value = 249.5
$ws.Range("K33").Value = 311.1875
$ws.Range("L33").Value = 249.5
$ws.Range("M33").Value = -82.1875
$ws.Range("N33").Value = -707.5
$ws.Range("H76").Value = 3300
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3300
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H96").Value = 1931.2632
$ws.Range("I96").Value = 1599.6
$ws.Range("J96").Value = 2299.7778
$ws.Range("K96").Value = 4798.799999999999
$ws.Range("L96").Value = 6899.3334
$ws.Range("M96").Value = -3425.799999999999
$ws.Range("N96").Value = -9645.3334
$ws.Range("H137").Value = 1272.2456
$ws.Range("I137").Value = 1083.6177
$ws.Range("K137").Value = 3250.8531
$ws.Range("M137").Value = -700.8531000000003
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 3000
$ws.Range("I35").Value = 3000
$ws.Range("K35").Value = 3000
$ws.Range("M35").Value = -2594
$ws.Range("H122").Value = 333930.66
$ws.Range("I122").Value = 333930.66
$ws.Range("K122").Value = 1001791.98
$ws.Range("M122").Value = -999341.98
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -887
$ws.Range("H94").Value = 469.8889
$ws.Range("I94").Value = 426.33334
$ws.Range("J94").Value = 491.66666
$ws.Range("K94").Value = 426.33334
$ws.Range("L94").Value = 491.66666
$ws.Range("M94").Value = 24.66665999999998
$ws.Range("N94").Value = -1393.66666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4243.4165
$ws.Range("I31").Value = 1619.4333
$ws.Range("J31").Value = 5701.185
$ws.Range("K31").Value = 1619.4333
$ws.Range("L31").Value = 5701.185
$ws.Range("M31").Value = -1324.4333
$ws.Range("N31").Value = -6291.185
$ws.Range("H34").Value = 4243.4165
$ws.Range("I34").Value = 1619.4333
$ws.Range("J34").Value = 5701.185
$ws.Range("K34").Value = 1619.4333
$ws.Range("L34").Value = 5701.185
$ws.Range("M34").Value = -1417.4333
$ws.Range("N34").Value = -6105.185
$ws.Range("H41").Value = 13993.667
$ws.Range("J41").Value = 15612.4
$ws.Range("L41").Value = 15612.4
$ws.Range("N41").Value = -16468.4
$ws.Range("H51").Value = 16399.4
$ws.Range("J51").Value = 17999.25
$ws.Range("L51").Value = 17999.25
$ws.Range("N51").Value = -19471.25
$ws.Range("H59").Value = 19199.2
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 19199.2
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 19199.2
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -21489.2
$ws.Range("H60").Value = 14342.417
$ws.Range("I60").Value = 700
$ws.Range("K60").Value = 700
$ws.Range("M60").Value = -189
$ws.Range("H61").Value = 16399.4
$ws.Range("J61").Value = 17999.25
$ws.Range("L61").Value = 17999.25
$ws.Range("N61").Value = -18695.25
$ws.Range("H68").Value = 21210.766
$ws.Range("J68").Value = 21210.766
$ws.Range("L68").Value = 21210.766
$ws.Range("N68").Value = -22708.766
$ws.Range("H71").Value = 21210.766
$ws.Range("J71").Value = 21210.766
$ws.Range("L71").Value = 63632.298
$ws.Range("N71").Value = -71120.298
$ws.Range("H74").Value = 20323.889
$ws.Range("J74").Value = 20323.889
$ws.Range("L74").Value = 20323.889
$ws.Range("N74").Value = -22071.889
$ws.Range("H77").Value = 20323.889
$ws.Range("J77").Value = 20323.889
$ws.Range("L77").Value = 60971.667
$ws.Range("N77").Value = -69707.667
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1299.75
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 1599.6666
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 4798.9998
$ws.Range("M17").Value = -1031
$ws.Range("N17").Value = -5136.9998
$ws.Range("H20").Value = 965
$ws.Range("I20").Value = 300
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 900
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -673
$ws.Range("N20").Value = -3454
$ws.Range("H34").Value = 10204646
$ws.Range("I34").Value = 228.33333
$ws.Range("J34").Value = 10870152
$ws.Range("K34").Value = 684.99999
$ws.Range("L34").Value = 32610456
$ws.Range("M34").Value = -600.99999
$ws.Range("N34").Value = -32610624
$ws.Range("H39").Value = 2097.95
$ws.Range("J39").Value = 2097.95
$ws.Range("L39").Value = 6293.849999999999
$ws.Range("N39").Value = -6881.849999999999
$ws.Range("H55").Value = 1384.75
$ws.Range("I55").Value = 1133.3334
$ws.Range("J55").Value = 1429.1177
$ws.Range("K55").Value = 3400.0002
$ws.Range("L55").Value = 4287.3531
$ws.Range("M55").Value = -3223.0002
$ws.Range("N55").Value = -4641.3531
$ws.Range("H70").Value = 4272.5713
$ws.Range("I70").Value = 4272.5713
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 12817.7139
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -12502.7139
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 4272.5713
$ws.Range("I73").Value = 4272.5713
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 12817.7139
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -11725.7139
$ws.Range("N73").ClearContents()
$ws.Range("H121").Value = 2940.8333
$ws.Range("I121").Value = 349.5
$ws.Range("J121").Value = 4236.5
$ws.Range("K121").Value = 1048.5
$ws.Range("L121").Value = 12709.5
$ws.Range("M121").Value = 261.5
$ws.Range("N121").Value = -15329.5
$ws.Range("H122").Value = 2968.1086
$ws.Range("I122").Value = 445.25
$ws.Range("J122").Value = 4908.769
$ws.Range("K122").Value = 4007.25
$ws.Range("L122").Value = 44178.921
$ws.Range("M122").Value = -1557.25
$ws.Range("N122").Value = -49078.921
$ws.Range("H134").Value = 6809.2095
$ws.Range("I134").Value = 3797.6924
$ws.Range("J134").Value = 8114.2
$ws.Range("K134").Value = 11393.0772
$ws.Range("L134").Value = 24342.6
$ws.Range("M134").Value = -6323.0772
$ws.Range("N134").Value = -34482.6
$ws.Range("H137").Value = 41992.43
$ws.Range("I137").Value = 7987.5293
$ws.Range("J137").Value = 94545.45
$ws.Range("K137").Value = 23962.5879
$ws.Range("L137").Value = 283636.35
$ws.Range("M137").Value = -18862.5879
$ws.Range("N137").Value = -293836.35
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2351620
$ws.Range("I80").Value = 3002160
$ws.Range("J80").Value = 400000
$ws.Range("K80").Value = 3002160
$ws.Range("L80").Value = 400000
$ws.Range("M80").Value = -3001162
$ws.Range("N80").Value = -401996
$ws.Range("H83").Value = 2351620
$ws.Range("I83").Value = 3002160
$ws.Range("J83").Value = 400000
$ws.Range("K83").Value = 15010800
$ws.Range("L83").Value = 2000000
$ws.Range("M83").Value = -15005808
$ws.Range("N83").Value = -2009984
$ws.Range("H102").Value = 2022.55
$ws.Range("I102").Value = 1907.6666
$ws.Range("J102").Value = 2367.2
$ws.Range("K102").Value = 1907.6666
$ws.Range("L102").Value = 2367.2
$ws.Range("M102").Value = -285.6666
$ws.Range("N102").Value = -5611.2
$ws.Range("H122").Value = 7340.15
$ws.Range("I122").Value = 8866.933999999999
$ws.Range("J122").Value = 2759.8
$ws.Range("K122").Value = 26600.802
$ws.Range("L122").Value = 8279.400000000001
$ws.Range("M122").Value = -24150.802
$ws.Range("N122").Value = -13179.4
$ws.Range("H132").Value = 25643294
$ws.Range("I132").Value = 30304808
$ws.Range("J132").Value = 4976
$ws.Range("K132").Value = 90914424
$ws.Range("L132").Value = 14928
$ws.Range("M132").Value = -90911894
$ws.Range("N132").Value = -19988
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H82").Value = 2077.75
$ws.Range("I82").Value = 1700
$ws.Range("J82").Value = 2203.6667
$ws.Range("K82").Value = 1700
$ws.Range("L82").Value = 2203.6667
$ws.Range("M82").Value = -1339
$ws.Range("N82").Value = -2925.6667
$ws.Range("H85").Value = 2077.75
$ws.Range("I85").Value = 1700
$ws.Range("J85").Value = 2203.6667
$ws.Range("K85").Value = 1700
$ws.Range("L85").Value = 2203.6667
$ws.Range("M85").Value = -452
$ws.Range("N85").Value = -4699.6667
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 49776.25
$ws.Range("J70").Value = 49776.25
$ws.Range("L70").Value = 49776.25
$ws.Range("N70").Value = -50406.25
$ws.Range("H73").Value = 49776.25
$ws.Range("J73").Value = 49776.25
$ws.Range("L73").Value = 49776.25
$ws.Range("N73").Value = -51960.25
$ws.Range("H113").Value = 2188.889
$ws.Range("I113").Value = 2566.6667
$ws.Range("J113").Value = 1433.3334
$ws.Range("K113").Value = 7700.000100000001
$ws.Range("L113").Value = 4300.0002
$ws.Range("M113").Value = -5530.000100000001
$ws.Range("N113").Value = -8640.0002
$ws.Range("H125").Value = 59481.25
$ws.Range("J125").Value = 60837.145
$ws.Range("L125").Value = 60837.145
$ws.Range("N125").Value = -70677.14499999999
$ws.Range("H137").Value = 47855
$ws.Range("J137").Value = 47855
$ws.Range("L137").Value = 47855
$ws.Range("N137").Value = -58055
